# Daten aktualisiert am 2024-01-11
# Append three new ticker rows at the end of the existing data in column A.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A241").Value = "IMX-USD"
$ws.Range("A242").Value = "MNT-USD"
$ws.Range("A243").Value = "GRT-USD"
